$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The "Price" column (D) holds plain text values like "27.920.36" or "1.000"
# that Excel would otherwise auto-convert to a number (losing the multi-dot
# grouping / trailing zeros / exact formatting). Prefix each new value with a
# leading quote so it's entered as text (matching the existing string cells on
# the sheet), then reset the cell style back to "Normal" so the quote-prefix
# formatting doesn't leave a visible style change on the cell.
#
# The "Volume(1h)" column (E) values are plain padded percentage strings; they
# already survive a normal .Value assignment as text (the surrounding spaces
# keep Excel from treating them as a numeric percentage), so no extra handling
# is needed there.

$ws.Range("D2").Formula = "'27.920.36"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.10%  "
$ws.Range("D3").Formula = "'1.767.95"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.73%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").Formula = "'328.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.24%  "
$ws.Range("D6").Formula = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Formula = "'0.4539"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.74%  "
$ws.Range("D8").Formula = "'0.3530"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("D9").Formula = "'41.92"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.76%  "
$ws.Range("D10").Formula = "'0.07387"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.89%  "
$ws.Range("D11").Formula = "'1.095"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.21%  "
$ws.Range("D12").Formula = "'1.000"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.11%  "
$ws.Range("D13").Formula = "'20.72"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.02%  "
$ws.Range("D14").Formula = "'6.011"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.54%  "
$ws.Range("D15").Formula = "'7.184"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.59%  "
$ws.Range("D16").Formula = "'1.766.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.89%  "
$ws.Range("D17").Formula = "'92.60"
$ws.Range("D17").Style = "Normal"
$ws.Range("E18").Value = "  +0.43%  "
$ws.Range("D19").Formula = "'0.06436"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.32%  "
$ws.Range("E20").Value = "  -0.04%  "
$ws.Range("D21").Formula = "'16.94"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Formula = "'5.770"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.41%  "
$ws.Range("D23").Formula = "'27.960.31"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +1.10%  "
$ws.Range("D24").Formula = "'11.20"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.03%  "
$ws.Range("D25").Formula = "'2.096"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.26%  "
$ws.Range("D26").Formula = "'159.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.10%  "
$ws.Range("D27").Formula = "'20.15"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.03%  "
$ws.Range("D28").Formula = "'1.980.08"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.37%  "
$ws.Range("D29").Formula = "'2.161"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +3.87%  "
$ws.Range("D30").Formula = "'124.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.83%  "
$ws.Range("D31").Formula = "'1.076"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Formula = "'0.09193"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.16%  "
$ws.Range("D33").Formula = "'5.620"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("D34").Formula = "'3.665"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.24%  "
$ws.Range("D35").Formula = "'11.84"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +1.31%  "
$ws.Range("D36").Formula = "'0.02284"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.16%  "
$ws.Range("D37").Formula = "'0.06120"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.52%  "
$ws.Range("D38").Formula = "'0.2093"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Formula = "'4.950"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.44%  "
$ws.Range("D40").Formula = "'0.6258"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.43%  "
$ws.Range("D41").Formula = "'1.179"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.43%  "
$ws.Range("D42").Formula = "'1.381"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.75%  "
$ws.Range("D43").Formula = "'7.809"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.70%  "
$ws.Range("D44").Formula = "'13.18"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").Formula = "'3.735"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.59%  "
$ws.Range("D46").Formula = "'0.5846"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.18%  "
$ws.Range("D47").Formula = "'122.62"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.08%  "
$ws.Range("D48").Formula = "'1.935"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.02%  "
$ws.Range("D49").Formula = "'1.130"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.20%  "
$ws.Range("D50").Formula = "'0.06826"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.95%  "
$ws.Range("E51").Value = "  +2.02%  "
